$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.498.46"
$ws.Range("E2").Value = "'  -0.44%  "
$ws.Range("D3").Value = "'1.819.47"
$ws.Range("E3").Value = "'  -0.55%  "
$ws.Range("D5").Value = "'316.55"
$ws.Range("E5").Value = "'  +0.04%  "
$ws.Range("E6").Value = "'  +0.18%  "
$ws.Range("D7").Value = "'0.5158"
$ws.Range("E7").Value = "'  -3.47%  "
$ws.Range("D8").Value = "'0.3879"
$ws.Range("E8").Value = "'  -3.03%  "
$ws.Range("D9").Value = "'0.08456"
$ws.Range("E9").Value = "'  +8.73%  "
$ws.Range("D10").Value = "'41.85"
$ws.Range("E10").Value = "'  -0.45%  "
$ws.Range("E11").Value = "'  -1.05%  "
$ws.Range("D12").Value = "'6.427"
$ws.Range("E12").Value = "'  +1.44%  "
$ws.Range("D13").Value = "'21.00"
$ws.Range("E13").Value = "'  -1.18%  "
$ws.Range("D14").Value = "'1.003"
$ws.Range("E14").Value = "'  +0.16%  "
$ws.Range("D15").Value = "'7.500"
$ws.Range("E15").Value = "'  -1.33%  "
$ws.Range("D16").Value = "'1.818.47"
$ws.Range("E16").Value = "'  -0.35%  "
$ws.Range("D17").Value = "'0.00001134"
$ws.Range("E17").Value = "'  +3.65%  "
$ws.Range("D18").Value = "'92.80"
$ws.Range("E18").Value = "'  -0.46%  "
$ws.Range("E19").Value = "'  +1.37%  "
$ws.Range("D20").Value = "'17.72"
$ws.Range("E20").Value = "'  -0.62%  "
$ws.Range("E21").Value = "'  +0.15%  "
$ws.Range("D22").Value = "'6.079"
$ws.Range("E22").Value = "'  -0.48%  "
$ws.Range("D23").Value = "'28.541.93"
$ws.Range("E23").Value = "'  -0.31%  "
$ws.Range("D24").Value = "'11.38"
$ws.Range("E24").Value = "'  +1.36%  "
$ws.Range("D25").Value = "'2.275"
$ws.Range("E25").Value = "'  +1.93%  "
$ws.Range("D26").Value = "'21.00"
$ws.Range("E26").Value = "'  +0.62%  "
$ws.Range("D27").Value = "'159.22"
$ws.Range("E27").Value = "'  +1.63%  "
$ws.Range("D28").Value = "'2.031.08"
$ws.Range("E28").Value = "'  -0.27%  "
$ws.Range("D29").Value = "'2.412"
$ws.Range("E29").Value = "'  -0.70%  "
$ws.Range("D30").Value = "'125.87"
$ws.Range("E30").Value = "'  +0.10%  "
$ws.Range("D31").Value = "'0.1085"
$ws.Range("E31").Value = "'  -3.75%  "
$ws.Range("E32").Value = "'  -5.87%  "
$ws.Range("D33").Value = "'5.734"
$ws.Range("E33").Value = "'  -0.64%  "
$ws.Range("D34").Value = "'0.07482"
$ws.Range("E34").Value = "'  +1.39%  "
$ws.Range("D35").Value = "'3.679"
$ws.Range("E35").Value = "'  +0.46%  "
$ws.Range("D36").Value = "'0.2233"
$ws.Range("E36").Value = "'  -2.09%  "
$ws.Range("E37").Value = "'  +0.24%  "
$ws.Range("D38").Value = "'5.198"
$ws.Range("E38").Value = "'  -0.46%  "
$ws.Range("D39").Value = "'8.756"
$ws.Range("E39").Value = "'  -2.15%  "
$ws.Range("D40").Value = "'0.6321"
$ws.Range("E40").Value = "'  +0.06%  "
$ws.Range("E41").Value = "'  -1.55%  "
$ws.Range("D42").Value = "'1.196"
$ws.Range("E42").Value = "'  -0.29%  "
$ws.Range("D43").Value = "'1.401"
$ws.Range("E43").Value = "'  +0.53%  "
$ws.Range("E44").Value = "'  -0.29%  "
$ws.Range("D45").Value = "'3.778"
$ws.Range("E45").Value = "'  +1.72%  "
$ws.Range("D46").Value = "'0.5932"
$ws.Range("E46").Value = "'  -0.46%  "
$ws.Range("D47").Value = "'126.10"
$ws.Range("D48").Value = "'1.992"
$ws.Range("E48").Value = "'  -0.62%  "
$ws.Range("D49").Value = "'1.199"
$ws.Range("E49").Value = "'  +0.33%  "
$ws.Range("D50").Value = "'0.06977"
$ws.Range("E50").Value = "'  +0.07%  "
$ws.Range("D51").Value = "'74.34"
$ws.Range("E51").Value = "'  -0.49%  "
